$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the literal time values in A1/B1 (keep the time number format already applied)
$ws.Range("A1").Value = 0.25
$ws.Range("B1").Value = 0.5

# Update the data validations whose formulas use literal time constants
# (previously 1:00:00 / 2:30:00, i.e. 0.041666... / 0.104166...)
$ranges = @("A2:A10", "B2:B10", "C2:C10", "D2:D10", "E2:E10", "F2:F10", "G2:G10", "H2:H10")
foreach ($addr in $ranges) {
    $rng = $ws.Range($addr)
    $dv = $rng.Validation
    $type = $dv.Type
    $operator = $dv.Operator
    if ($addr -eq "G2:G10" -or $addr -eq "H2:H10") {
        $dv.Modify($type, $dv.AlertStyle, $operator, "0.25", "0.5")
    } else {
        $dv.Modify($type, $dv.AlertStyle, $operator, "0.25")
    }
}
